# Apply updated "current" counts and recompute "change" values
# for the agency count comparison sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new "current" (column B) value.
# "previous" (column C) stays the same; "change" (column D) = B - C.
$updates = @{
    6  = 29   # Department of Commerce
    7  = 40   # Department of Defense
    8  = 27   # Department of Education
    10 = 51   # Department of Health and Human Services
    14 = 30   # Department of Labor
    15 = 34   # Department of State
    21 = 11   # Executive Office of the President, Management and Administration
    26 = 28   # Intelligence Community
    35 = 20   # Small Business Administration
    38 = 23   # United States Department of Agriculture
    39 = 3    # United States Digital Service
}

foreach ($row in $updates.Keys) {
    $newCurrent = $updates[$row]
    $previous = $ws.Cells.Item($row, 3).Value2

    $ws.Cells.Item($row, 2).Value = $newCurrent
    $ws.Cells.Item($row, 4).Value = $newCurrent - $previous
}
